$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The recomputed (new TPM) results no longer contain the MuSCs -> ECs
# sending/target pairing that used to be row 2, so drop it; the remaining
# rows (old rows 3 and 4) shift up to become rows 2 and 3.
$ws.Rows("2:2").Delete()

# With the new TPM data the sending cluster for the remaining pairs is now
# "ECs" (it used to be "MuSCs"); the target clusters (FAPs, then MuSCs) stay
# the same as before for these two rows.
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"

# Row 2 (Target cluster = FAPs): refresh the ligand-level stats (shared
# across rows, now reflecting the ECs sending cluster) and the
# receptor/edge-derived specificity stats with the new TPM-based numbers.
$ws.Range("G2").Value = 0.072366
$ws.Range("H2").Value = 0.217098
$ws.Range("M2").Value = [double]"0.0005823333333333334"
$ws.Range("N2").Value = 0.001747
$ws.Range("O2").Value = [double]"9.625315715314126E-05"
$ws.Range("P2").Value = [double]"9.625315715314125E-05"
$ws.Range("Q2").Value = [double]"4.2141134E-05"
$ws.Range("R2").Value = [double]"0.0003792702060000001"
$ws.Range("S2").Value = [double]"9.625315715314126E-05"
$ws.Range("T2").Value = [double]"9.625315715314125E-05"

# Row 3 (Target cluster = MuSCs): refresh the ligand-level stats and all
# receptor/edge-derived stats with the new TPM-based numbers.
$ws.Range("G3").Value = 0.072366
$ws.Range("H3").Value = 0.217098
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.049435666666667
$ws.Range("N3").Value = 18.148307
$ws.Range("O3").Value = [double]"0.9999037468428469"
$ws.Range("P3").Value = [double]"0.9999037468428468"
$ws.Range("Q3").Value = [double]"0.437773461454"
$ws.Range("R3").Value = [double]"3.939961153086"
$ws.Range("S3").Value = [double]"0.9999037468428469"
$ws.Range("T3").Value = [double]"0.9999037468428468"
